# Applies the cryptos.xlsx data refresh (GitHub Actions "Updated cryptos list" commit).
# Column D ("Price") and E ("Volume(1h)") values are refreshed for most rows, and two
# row pairs swap order: Avalanche/TRON (rows 14-15) and InternetComputer(DFINITY)/
# Fetch.AI (rows 31-32).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell reference -> new value to write.
$updates = [ordered]@{
    'D2' = '60.780.01'
    'E2' = '  -1.74%  '
    'D3' = '3.384.61'
    'E3' = '  -2.20%  '
    'D4' = '0.999'
    'E4' = '  +0.00%  '
    'D5' = '571.83'
    'D6' = '141.76'
    'E6' = '  -4.87%  '
    'E7' = '  +0.06%  '
    'D8' = '3.384.13'
    'E8' = '  -2.29%  '
    'D9' = '0.474'
    'E9' = '  -0.18%  '
    'D10' = '7.51'
    'E10' = '  -3.40%  '
    'D11' = '0.124'
    'E11' = '  -1.05%  '
    'E12' = '  +0.41%  '
    'D13' = '3.963.39'
    'E13' = '  -2.17%  '
    'B14' = 'Avalanche'
    'C14' = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
    'D14' = '28.10'
    'E14' = '  -0.65%  '
    'B15' = 'TRON'
    'C15' = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
    'D15' = '0.124'
    'E15' = '  +0.94%  '
    'E16' = '  -2.99%  '
    'D17' = '3.384.35'
    'E17' = '  -2.13%  '
    'D18' = '60.895.62'
    'E18' = '  -1.58%  '
    'D19' = '6.30'
    'E19' = '  -0.97%  '
    'D20' = '14.14'
    'E20' = '  -1.65%  '
    'D21' = '8.96'
    'E21' = '  -5.34%  '
    'D22' = '388.61'
    'E22' = '  +0.74%  '
    'D23' = '0.561'
    'E23' = '  -1.66%  '
    'D24' = '73.46'
    'E24' = '  +1.06%  '
    'E25' = '  +0.21%  '
    'E26' = '  -4.07%  '
    'D27' = '3.523.18'
    'E27' = '  -2.08%  '
    'E28' = '  -1.54%  '
    'E29' = '  -0.15%  '
    'D30' = '7.39'
    'E30' = '  -5.87%  '
    'B31' = 'Fetch.AI'
    'C31' = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
    'D31' = '1.46'
    'E31' = '  -4.51%  '
    'B32' = 'InternetComputer(DFINITY)'
    'C32' = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
    'D32' = '8.09'
    'E32' = '  -2.06%  '
    'D33' = '2.17'
    'E33' = '  -0.17%  '
    'E34' = '  -0.02%  '
    'D35' = '23.82'
    'E35' = '  -0.67%  '
    'D36' = '6.93'
    'E36' = '  -1.99%  '
    'D37' = '3.412.41'
    'E37' = '  -1.94%  '
    'D38' = '166.89'
    'E38' = '  +0.35%  '
    'D39' = '5.05'
    'E39' = '  -3.58%  '
    'E40' = '  -3.16%  '
    'D41' = '0.0780'
    'D42' = '26.80'
    'E42' = '  +3.44%  '
    'D43' = '0.784'
    'E43' = '  -1.61%  '
    'E44' = '  -0.01%  '
    'D45' = '4.46'
    'E45' = '  -0.62%  '
    'D46' = '41.72'
    'E46' = '  -1.53%  '
    'E47' = '  -2.48%  '
    'D48' = '2.552.53'
    'E48' = '  -2.03%  '
    'D49' = '1.13'
    'E49' = '  -4.30%  '
    'D50' = '6.83'
    'E50' = '  -2.06%  '
    'D51' = '22.92'
    'E51' = '  -1.91%  '
}

foreach ($cellRef in $updates.Keys) {
    $newValue = $updates[$cellRef]

    # Several "Price" values are plain numeric-looking strings (e.g. "28.10",
    # "0.0780"). Assigning them straight to .Value would make Excel coerce them
    # into real numbers and lose formatting such as trailing zeros. Prefix those
    # with a single quote (exactly like typing '28.10 into a cell) so they are
    # stored as text, matching the rest of the column. Values with more than one
    # "." (e.g. "60.780.01") are never valid numbers, so Excel already keeps
    # those as plain text without any extra help.
    if ($newValue -match '^[+-]?\d+(\.\d+)?$') {
        $ws.Range($cellRef).Value = "'" + $newValue
    } else {
        $ws.Range($cellRef).Value = $newValue
    }
}
